$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 62223840
$ws.Range("I70").Value = 74668110
$ws.Range("J70").Value = 2500
$ws.Range("K70").Value = 224004330
$ws.Range("L70").Value = 7500
$ws.Range("M70").Value = -224004060
$ws.Range("N70").Value = -8040

$ws.Range("H73").Value = 62223840
$ws.Range("I73").Value = 74668110
$ws.Range("J73").Value = 2500
$ws.Range("K73").Value = 224004330
$ws.Range("L73").Value = 7500
$ws.Range("M73").Value = -224003394
$ws.Range("N73").Value = -9372

$ws.Range("H127").Value = 2095.1428
$ws.Range("I127").Value = 975.75
$ws.Range("J127").Value = 2784
$ws.Range("K127").Value = 2927.25
$ws.Range("L127").Value = 8352
$ws.Range("M127").Value = 2032.75
$ws.Range("N127").Value = -18272

$ws.Range("H137").Value = 272046.7
$ws.Range("I137").Value = 518697.22
$ws.Range("J137").Value = 1905.619
$ws.Range("K137").Value = 1556091.66
$ws.Range("L137").Value = 5716.857
$ws.Range("M137").Value = -1553541.66
$ws.Range("N137").Value = -10816.857

$ws.Range("H138").Value = 2086.7532
$ws.Range("I138").Value = 908.7308
$ws.Range("J138").Value = 2643.6365
$ws.Range("K138").Value = 2726.1924
$ws.Range("L138").Value = 7930.9095
$ws.Range("M138").Value = 2413.8076
$ws.Range("N138").Value = -18210.9095

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1048.6923
$ws.Range("I45").Value = 1124.6
$ws.Range("J45").Value = 1001.25
$ws.Range("K45").Value = 1124.6
$ws.Range("L45").Value = 1001.25
$ws.Range("M45").Value = -747.5999999999999
$ws.Range("N45").Value = -1755.25

$ws.Range("H97").Value = 586.86957
$ws.Range("I97").Value = 477.18182
$ws.Range("J97").Value = 3000
$ws.Range("K97").Value = 477.18182
$ws.Range("L97").Value = 3000
$ws.Range("M97").Value = 18.81817999999998
$ws.Range("N97").Value = -3992

$ws.Range("H110").Value = 851.0333000000001
$ws.Range("I110").Value = 794.4706
$ws.Range("J110").Value = 925
$ws.Range("K110").Value = 794.4706
$ws.Range("L110").Value = 925
$ws.Range("M110").Value = 1250.5294
$ws.Range("N110").Value = -5015

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4812.0938
$ws.Range("I20").Value = 6044.909
$ws.Range("J20").Value = 2099.9
$ws.Range("K20").Value = 6044.909
$ws.Range("L20").Value = 2099.9
$ws.Range("M20").Value = -5797.909
$ws.Range("N20").Value = -2593.9

$ws.Range("H86").Value = 1632.4
$ws.Range("I86").Value = 1627.5714
$ws.Range("J86").Value = 1700
$ws.Range("K86").Value = 1627.5714
$ws.Range("L86").Value = 1700
$ws.Range("M86").Value = -504.5714
$ws.Range("N86").Value = -3946

$ws.Range("H89").Value = 1632.4
$ws.Range("I89").Value = 1627.5714
$ws.Range("J89").Value = 1700
$ws.Range("K89").Value = 8137.857
$ws.Range("L89").Value = 8500
$ws.Range("M89").Value = -2521.857
$ws.Range("N89").Value = -19732

$ws.Range("H94").Value = 438.125
$ws.Range("I94").Value = 376.25
$ws.Range("J94").Value = 747.5
$ws.Range("K94").Value = 376.25
$ws.Range("L94").Value = 747.5
$ws.Range("M94").Value = 74.75
$ws.Range("N94").Value = -1649.5

$ws.Range("H99").Value = 932.2222
$ws.Range("I99").Value = 898.3333
$ws.Range("J99").Value = 1000
$ws.Range("K99").Value = 898.3333
$ws.Range("L99").Value = 1000
$ws.Range("M99").Value = 599.6667
$ws.Range("N99").Value = -3996

$ws.Range("H105").Value = 2000
$ws.Range("I105").Value = 2000
$ws.Range("J105").Value = 2000
$ws.Range("K105").Value = 2000
$ws.Range("L105").Value = 2000
$ws.Range("M105").Value = -253
$ws.Range("N105").Value = -5494

$ws.Range("H107").Value = 71429500
$ws.Range("I107").Value = 100001050
$ws.Range("J107").Value = 642
$ws.Range("K107").Value = 100001050
$ws.Range("L107").Value = 642
$ws.Range("M107").Value = -99999130
$ws.Range("N107").Value = -4482

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8629.764999999999
$ws.Range("I31").Value = 3110.7827
$ws.Range("J31").Value = 20169.455
$ws.Range("K31").Value = 3110.7827
$ws.Range("L31").Value = 20169.455
$ws.Range("M31").Value = -2815.7827
$ws.Range("N31").Value = -20759.455

$ws.Range("H34").Value = 8629.764999999999
$ws.Range("I34").Value = 3110.7827
$ws.Range("J34").Value = 20169.455
$ws.Range("K34").Value = 3110.7827
$ws.Range("L34").Value = 20169.455
$ws.Range("M34").Value = -2908.7827
$ws.Range("N34").Value = -20573.455

$ws.Range("H58").Value = 2362433.5
$ws.Range("I58").Value = 3425882.2
$ws.Range("J58").Value = 11651.947
$ws.Range("K58").Value = 3425882.2
$ws.Range("L58").Value = 11651.947
$ws.Range("M58").Value = -3425679.2
$ws.Range("N58").Value = -12057.947

$ws.Range("H132").Value = 10106421
$ws.Range("I132").Value = 15875370
$ws.Range("J132").Value = 10760.333
$ws.Range("K132").Value = 47626110
$ws.Range("L132").Value = 32280.999
$ws.Range("M132").Value = -47623580
$ws.Range("N132").Value = -37340.999

$ws.Range("H136").Value = 2362433.5
$ws.Range("I136").Value = 3425882.2
$ws.Range("J136").Value = 11651.947
$ws.Range("K136").Value = 10277646.6
$ws.Range("L136").Value = 34955.841
$ws.Range("M136").Value = -10275096.6
$ws.Range("N136").Value = -40055.841

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 288.75
$ws.Range("I7").Value = 227
$ws.Range("J7").Value = 332.85715
$ws.Range("K7").Value = 681
$ws.Range("L7").Value = 998.5714499999999
$ws.Range("M7").Value = -569
$ws.Range("N7").Value = -1222.57145

$ws.Range("H107").Value = 6521.1055
$ws.Range("I107").Value = 249.625
$ws.Range("J107").Value = 11082.182
$ws.Range("K107").Value = 748.875
$ws.Range("L107").Value = 33246.546
$ws.Range("M107").Value = 1171.125
$ws.Range("N107").Value = -37086.546

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H59").Value = 80110
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 80110
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 80110
$ws.Range("N59").Value = -81276

$ws.Range("H70").Value = 37819.633
$ws.Range("I70").Value = 49836.316
$ws.Range("J70").Value = 4773.75
$ws.Range("K70").Value = 49836.316
$ws.Range("L70").Value = 4773.75
$ws.Range("M70").Value = -49566.316
$ws.Range("N70").Value = -5313.75

$ws.Range("H73").Value = 37819.633
$ws.Range("I73").Value = 49836.316
$ws.Range("J73").Value = 4773.75
$ws.Range("K73").Value = 49836.316
$ws.Range("L73").Value = 4773.75
$ws.Range("M73").Value = -48900.316
$ws.Range("N73").Value = -6645.75

$ws.Range("H80").Value = 2711.7646
$ws.Range("I80").Value = 2233.3333
$ws.Range("J80").Value = 2758.0645
$ws.Range("K80").Value = 2233.3333
$ws.Range("L80").Value = 2758.0645
$ws.Range("M80").Value = -1235.3333

$ws.Range("H83").Value = 2711.7646
$ws.Range("I83").Value = 2233.3333
$ws.Range("J83").Value = 2758.0645
$ws.Range("K83").Value = 11166.6665
$ws.Range("L83").Value = 13790.3225
$ws.Range("M83").Value = -6174.666499999999

$ws.Range("H97").Value = 153846820
$ws.Range("I97").Value = 166667360
$ws.Range("J97").Value = 142857800
$ws.Range("K97").Value = 166667360
$ws.Range("L97").Value = 142857800
$ws.Range("M97").Value = -166666864
$ws.Range("N97").Value = -142858792

$ws.Range("H113").Value = 1719.238
$ws.Range("I113").Value = 1338.5385
$ws.Range("J113").Value = 2337.875
$ws.Range("K113").Value = 1338.5385
$ws.Range("L113").Value = 2337.875
$ws.Range("M113").Value = 831.4614999999999
$ws.Range("N113").Value = -6677.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5510
$ws.Range("I46").Value = 1033.3334
$ws.Range("J46").Value = 15102.857
$ws.Range("K46").Value = 1033.3334
$ws.Range("L46").Value = 15102.857
$ws.Range("M46").Value = -845.3334
$ws.Range("N46").Value = -15478.857

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 30000
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 30000
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 30000
$ws.Range("N54").Value = -31040

$ws.Range("H62").Value = 6000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 6000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 6000
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -7248

$ws.Range("H65").Value = 6000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 6000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 30000
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -36240

$ws.Range("H132").Value = 17852906
$ws.Range("I132").Value = 6251675
$ws.Range("J132").Value = 64257830
$ws.Range("K132").Value = 18755025
$ws.Range("L132").Value = 192773490
$ws.Range("M132").Value = -18752495
$ws.Range("N132").Value = -192778550

$ws.Range("H136").Value = 11132995
$ws.Range("I136").Value = 5412342.5
$ws.Range("J136").Value = 45456908
$ws.Range("K136").Value = 16237027.5
$ws.Range("L136").Value = 136370724
$ws.Range("M136").Value = -16234477.5
$ws.Range("N136").Value = -136375824
